$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 226.83333
$ws.Range("I11").Value = 226.83333
$ws.Range("K11").Value = 226.83333
$ws.Range("M11").Value = -86.83332999999999
$ws.Range("H18").Value = 3405.3333
$ws.Range("I18").Value = 3405.3333
$ws.Range("K18").Value = 3405.3333
$ws.Range("M18").Value = -3121.3333
$ws.Range("H33").Value = 309.1111
$ws.Range("I33").Value = 332.875
$ws.Range("K33").Value = 332.875
$ws.Range("M33").Value = -103.875
$ws.Range("H40").Value = 29413842
$ws.Range("J40").Value = 71430750
$ws.Range("L40").Value = 71430750
$ws.Range("N40").Value = -71431100
$ws.Range("H64").Value = 14708384
$ws.Range("I64").Value = 29411764
$ws.Range("J64").Value = 5003
$ws.Range("K64").Value = 29411764
$ws.Range("L64").Value = 5003
$ws.Range("M64").Value = -29411516
$ws.Range("N64").Value = -5499
$ws.Range("H67").Value = 14708384
$ws.Range("I67").Value = 29411764
$ws.Range("J67").Value = 5003
$ws.Range("K67").Value = 29411764
$ws.Range("L67").Value = 5003
$ws.Range("M67").Value = -29410906
$ws.Range("N67").Value = -6719
$ws.Range("H80").Value = 6724729
$ws.Range("J80").Value = 7938907.5
$ws.Range("L80").Value = 23816722.5
$ws.Range("N80").Value = -23818718.5
$ws.Range("H83").Value = 6724729
$ws.Range("J83").Value = 7938907.5
$ws.Range("L83").Value = 71450167.5
$ws.Range("N83").Value = -71460151.5
$ws.Range("H96").Value = 1116716.2
$ws.Range("J96").Value = 2418000.5
$ws.Range("L96").Value = 7254001.5
$ws.Range("N96").Value = -7256747.5
$ws.Range("H97").Value = 10999
$ws.Range("J97").Value = 10999
$ws.Range("L97").Value = 32997
$ws.Range("N97").Value = -33989
$ws.Range("H98").Value = 3907695.8
$ws.Range("I98").Value = 4167875.2
$ws.Range("K98").Value = 4167875.2
$ws.Range("M98").Value = -4166377.2
$ws.Range("H112").Value = 2783483.8
$ws.Range("J112").Value = 4174645.8
$ws.Range("L112").Value = 12523937.4
$ws.Range("N112").Value = -12526153.4
$ws.Range("H113").Value = 4002.5
$ws.Range("I113").Value = 4002.5
$ws.Range("K113").Value = 4002.5
$ws.Range("M113").Value = -748.5
$ws.Range("H116").Value = 9192.679
$ws.Range("I116").Value = 6076.25
$ws.Range("J116").Value = 9712.083000000001
$ws.Range("K116").Value = 6076.25
$ws.Range("L116").Value = 9712.083000000001
$ws.Range("M116").Value = -2634.25
$ws.Range("N116").Value = -16596.083
$ws.Range("H122").Value = 3907695.8
$ws.Range("I122").Value = 4167875.2
$ws.Range("K122").Value = 12503625.6
$ws.Range("M122").Value = -12501175.6
$ws.Range("H125").Value = 2102.1428
$ws.Range("J125").Value = 2203
$ws.Range("L125").Value = 19827
$ws.Range("N125").Value = -24747
$ws.Range("H135").Value = 3361.0833
$ws.Range("I135").Value = 2410.3428
$ws.Range("J135").Value = 5920.769
$ws.Range("K135").Value = 21693.0852
$ws.Range("L135").Value = 53286.921
$ws.Range("M135").Value = -19158.0852
$ws.Range("N135").Value = -58356.921
$ws.Range("H137").Value = 3007
$ws.Range("I137").Value = 2846.375
$ws.Range("J137").Value = 3521
$ws.Range("K137").Value = 8539.125
$ws.Range("L137").Value = 10563
$ws.Range("M137").Value = -5989.125
$ws.Range("N137").Value = -15663
$ws.Range("H138").Value = 6692.273
$ws.Range("I138").Value = 3995.1875
$ws.Range("K138").Value = 11985.5625
$ws.Range("M138").Value = -6845.5625
$ws.Range("H141").Value = 13161994
$ws.Range("I141").Value = 15154922
$ws.Range("J141").Value = 8670.6
$ws.Range("K141").Value = 45464766
$ws.Range("L141").Value = 26011.8
$ws.Range("M141").Value = -45459586
$ws.Range("N141").Value = -36371.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 79989
$ws.Range("J30").Value = 79989
$ws.Range("L30").Value = 79989
$ws.Range("N30").Value = -80289
$ws.Range("H32").Value = 2723.1973
$ws.Range("I32").Value = 2453.8262
$ws.Range("K32").Value = 2453.8262
$ws.Range("M32").Value = -2166.8262
$ws.Range("H43").Value = 21442.889
$ws.Range("J43").Value = 24177.428
$ws.Range("L43").Value = 24177.428
$ws.Range("N43").Value = -24803.428
$ws.Range("H45").Value = 2663.6667
$ws.Range("I45").Value = 1687.25
$ws.Range("J45").Value = 3965.5557
$ws.Range("K45").Value = 1687.25
$ws.Range("L45").Value = 3965.5557
$ws.Range("M45").Value = -1310.25
$ws.Range("N45").Value = -4719.5557
$ws.Range("H61").Value = 10002381
$ws.Range("I61").Value = 13335775
$ws.Range("J61").Value = 2859394.8
$ws.Range("K61").Value = 13335775
$ws.Range("L61").Value = 2859394.8
$ws.Range("M61").Value = -13335563
$ws.Range("N61").Value = -2859818.8
$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -314
$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("M66").Value = -1568
$ws.Range("H74").Value = 2237.2
$ws.Range("I74").Value = 2051.7778
$ws.Range("J74").Value = 2714
$ws.Range("K74").Value = 2051.7778
$ws.Range("L74").Value = 2714
$ws.Range("M74").Value = -1177.7778
$ws.Range("N74").Value = -4462
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H77").Value = 2237.2
$ws.Range("I77").Value = 2051.7778
$ws.Range("J77").Value = 2714
$ws.Range("K77").Value = 10258.889
$ws.Range("L77").Value = 13570
$ws.Range("M77").Value = -5890.888999999999
$ws.Range("N77").Value = -22306
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H102").Value = 2313.84
$ws.Range("I102").Value = 1598.1666
$ws.Range("J102").Value = 4154.143
$ws.Range("K102").Value = 1598.1666
$ws.Range("L102").Value = 4154.143
$ws.Range("M102").Value = 23.83339999999998
$ws.Range("N102").Value = -7398.143
$ws.Range("H132").Value = 13640150
$ws.Range("I132").Value = 3747.4
$ws.Range("J132").Value = 42861012
$ws.Range("K132").Value = 11242.2
$ws.Range("L132").Value = 128583036
$ws.Range("M132").Value = -8712.200000000001
$ws.Range("N132").Value = -128588096
$ws.Range("H136").Value = 10002381
$ws.Range("I136").Value = 13335775
$ws.Range("J136").Value = 2859394.8
$ws.Range("K136").Value = 40007325
$ws.Range("L136").Value = 8578184.399999999
$ws.Range("M136").Value = -40004775
$ws.Range("N136").Value = -8583284.399999999
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 67999.89
$ws.Range("J57").Value = 67999.89
$ws.Range("L57").Value = 67999.89
$ws.Range("N57").Value = -69439.89
$ws.Range("H86").Value = 1104478.6
$ws.Range("I86").Value = 2752372.8
$ws.Range("K86").Value = 2752372.8
$ws.Range("M86").Value = -2751249.8
$ws.Range("H89").Value = 1104478.6
$ws.Range("I89").Value = 2752372.8
$ws.Range("K89").Value = 13761864
$ws.Range("M89").Value = -13756248
$ws.Range("H96").Value = 23332
$ws.Range("I96").Value = 23332
$ws.Range("K96").Value = 23332
$ws.Range("M96").Value = -20586
$ws.Range("H134").Value = 3335060.2
$ws.Range("I134").Value = 1615
$ws.Range("J134").Value = 33336066
$ws.Range("K134").Value = 4845
$ws.Range("L134").Value = 100008198
$ws.Range("M134").Value = -2310
$ws.Range("N134").Value = -100013268
$ws.Range("H136").Value = 67999.89
$ws.Range("J136").Value = 67999.89
$ws.Range("L136").Value = 67999.89
$ws.Range("N136").Value = -78199.89
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H139").Value = 70000
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 100000
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 100000
$ws.Range("M139").Value = -34860
$ws.Range("N139").Value = -110280
$ws.Range("H140").Value = 231999.25
$ws.Range("J140").Value = 231999.25
$ws.Range("L140").Value = 231999.25
$ws.Range("N140").Value = -242359.25
$ws.Range("H141").Value = 111799.8
$ws.Range("J141").Value = 111799.8
$ws.Range("L141").Value = 111799.8
$ws.Range("N141").Value = -122159.8
$ws.Range("N137").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 116.1
$ws.Range("J7").Value = 269.8
$ws.Range("L7").Value = 269.8
$ws.Range("N7").Value = -495.8
$ws.Range("H31").Value = 15887071
$ws.Range("I31").Value = 29435866
$ws.Range("J31").Value = 2277
$ws.Range("K31").Value = 29435866
$ws.Range("L31").Value = 2277
$ws.Range("M31").Value = -29435571
$ws.Range("N31").Value = -2867
$ws.Range("H34").Value = 15887071
$ws.Range("I34").Value = 29435866
$ws.Range("J34").Value = 2277
$ws.Range("K34").Value = 29435866
$ws.Range("L34").Value = 2277
$ws.Range("M34").Value = -29435664
$ws.Range("N34").Value = -2681
$ws.Range("H51").Value = 22290
$ws.Range("I51").Value = 22290
$ws.Range("K51").Value = 22290
$ws.Range("M51").Value = -21554
$ws.Range("H60").Value = 2041
$ws.Range("I60").Value = 2041
$ws.Range("K60").Value = 2041
$ws.Range("M60").Value = -1530
$ws.Range("H61").Value = 22290
$ws.Range("I61").Value = 22290
$ws.Range("K61").Value = 22290
$ws.Range("M61").Value = -21942
$ws.Range("H62").Value = 83339830
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 111116450
$ws.Range("K62").Value = 9999
$ws.Range("L62").Value = 111116450
$ws.Range("M62").Value = -9375
$ws.Range("N62").Value = -111117698
$ws.Range("H64").Value = 97399
$ws.Range("J64").Value = 97399
$ws.Range("L64").Value = 97399
$ws.Range("N64").Value = -97895
$ws.Range("H65").Value = 83339830
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 111116450
$ws.Range("K65").Value = 49995
$ws.Range("L65").Value = 555582250
$ws.Range("M65").Value = -46875
$ws.Range("N65").Value = -555588490
$ws.Range("H67").Value = 97399
$ws.Range("J67").Value = 97399
$ws.Range("L67").Value = 97399
$ws.Range("N67").Value = -99115
$ws.Range("H105").Value = 2405.7
$ws.Range("I105").Value = 2153.8572
$ws.Range("K105").Value = 2153.8572
$ws.Range("M105").Value = -406.8571999999999
$ws.Range("H132").Value = 2086.48
$ws.Range("I132").Value = 1872.409
$ws.Range("J132").Value = 3656.3333
$ws.Range("K132").Value = 5617.227000000001
$ws.Range("L132").Value = 10968.9999
$ws.Range("M132").Value = -3087.227000000001
$ws.Range("N132").Value = -16028.9999
$ws.Range("H134").Value = 3131.361
$ws.Range("I134").Value = 3131.361
$ws.Range("K134").Value = 9394.082999999999
$ws.Range("M134").Value = -6859.082999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 312.375
$ws.Range("I2").Value = 102.833336
$ws.Range("J2").Value = 438.1
$ws.Range("K2").Value = 617.000016
$ws.Range("L2").Value = 2628.6
$ws.Range("M2").Value = -504.000016
$ws.Range("N2").Value = -2854.6
$ws.Range("H3").Value = 9534.727999999999
$ws.Range("I3").Value = 7154.9
$ws.Range("K3").Value = 21464.7
$ws.Range("M3").Value = -21352.7
$ws.Range("H29").Value = 10165.875
$ws.Range("I29").Value = 497.5
$ws.Range("J29").Value = 13388.667
$ws.Range("K29").Value = 1492.5
$ws.Range("L29").Value = 40166.001
$ws.Range("M29").Value = -1215.5
$ws.Range("N29").Value = -40720.001
$ws.Range("H56").Value = 15245.724
$ws.Range("I56").Value = 15245.724
$ws.Range("K56").Value = 15245.724
$ws.Range("M56").Value = -14715.724
$ws.Range("H68").Value = 1891
$ws.Range("J68").Value = 2745.6667
$ws.Range("L68").Value = 8237.000100000001
$ws.Range("N68").Value = -9859.000100000001
$ws.Range("H71").Value = 1891
$ws.Range("J71").Value = 2745.6667
$ws.Range("L71").Value = 24711.0003
$ws.Range("N71").Value = -32823.0003
$ws.Range("H76").Value = 16448.666
$ws.Range("I76").Value = 8006.5
$ws.Range("J76").Value = 33333
$ws.Range("K76").Value = 24019.5
$ws.Range("L76").Value = 99999
$ws.Range("M76").Value = -23636.5
$ws.Range("N76").Value = -100765
$ws.Range("H79").Value = 16448.666
$ws.Range("I79").Value = 8006.5
$ws.Range("J79").Value = 33333
$ws.Range("K79").Value = 24019.5
$ws.Range("L79").Value = 99999
$ws.Range("M79").Value = -22693.5
$ws.Range("N79").Value = -102651
$ws.Range("H97").Value = 1361.75
$ws.Range("J97").Value = 1167.5454
$ws.Range("L97").Value = 3502.6362
$ws.Range("N97").Value = -4494.6362
$ws.Range("H104").Value = 20665.5
$ws.Range("I104").Value = 7998
$ws.Range("J104").Value = 33333
$ws.Range("K104").Value = 23994
$ws.Range("L104").Value = 99999
$ws.Range("M104").Value = -21373
$ws.Range("N104").Value = -105241
$ws.Range("H107").Value = 4794055
$ws.Range("I107").Value = 2210.7144
$ws.Range("J107").Value = 7589298
$ws.Range("K107").Value = 6632.1432
$ws.Range("L107").Value = 22767894
$ws.Range("M107").Value = -4712.1432
$ws.Range("N107").Value = -22771734
$ws.Range("H117").Value = 4557.933
$ws.Range("I117").Value = 233.33333
$ws.Range("J117").Value = 5639.0835
$ws.Range("K117").Value = 699.99999
$ws.Range("L117").Value = 16917.2505
$ws.Range("M117").Value = 2742.00001
$ws.Range("N117").Value = -23801.2505
$ws.Range("H122").Value = 41305
$ws.Range("I122").Value = 83437
$ws.Range("J122").Value = 7599.4
$ws.Range("K122").Value = 750933
$ws.Range("L122").Value = 68394.59999999999
$ws.Range("M122").Value = -748483
$ws.Range("N122").Value = -73294.59999999999
$ws.Range("H131").Value = 3431.742
$ws.Range("I131").Value = 2217.2856
$ws.Range("J131").Value = 5982.1
$ws.Range("K131").Value = 6651.8568
$ws.Range("L131").Value = 17946.3
$ws.Range("M131").Value = -1611.8568
$ws.Range("N131").Value = -28026.3
$ws.Range("H132").Value = 1817.875
$ws.Range("I132").Value = 1750.75
$ws.Range("J132").Value = 1885
$ws.Range("K132").Value = 15756.75
$ws.Range("L132").Value = 16965
$ws.Range("M132").Value = -13226.75
$ws.Range("N132").Value = -22025
$ws.Range("H134").Value = 9528.317999999999
$ws.Range("I134").Value = 2606.2354
$ws.Range("K134").Value = 7818.706200000001
$ws.Range("M134").Value = -2748.706200000001
$ws.Range("H136").Value = 7951.75
$ws.Range("I136").Value = 2870.1428
$ws.Range("K136").Value = 8610.428400000001
$ws.Range("M136").Value = -3510.428400000001
$ws.Range("H138").Value = 16380.75
$ws.Range("I138").Value = 17605.143
$ws.Range("K138").Value = 52815.429
$ws.Range("M138").Value = -47675.429
$ws.Range("H139").Value = 1923.3956
$ws.Range("I139").Value = 1876
$ws.Range("K139").Value = 5628
$ws.Range("M139").Value = -488

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 46.2
$ws.Range("J2").Value = 22.5
$ws.Range("L2").Value = 22.5
$ws.Range("N2").Value = -248.5
$ws.Range("H80").Value = 2305.5
$ws.Range("I80").Value = 1605
$ws.Range("J80").Value = 3006
$ws.Range("K80").Value = 1605
$ws.Range("L80").Value = 3006
$ws.Range("M80").Value = -607
$ws.Range("N80").Value = -5002
$ws.Range("H82").Value = 59999
$ws.Range("J82").Value = 59999
$ws.Range("L82").Value = 59999
$ws.Range("N82").Value = -60765
$ws.Range("H83").Value = 2305.5
$ws.Range("I83").Value = 1605
$ws.Range("J83").Value = 3006
$ws.Range("K83").Value = 8025
$ws.Range("L83").Value = 15030
$ws.Range("M83").Value = -3033
$ws.Range("N83").Value = -25014
$ws.Range("H85").Value = 59999
$ws.Range("J85").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("N85").Value = -62651
$ws.Range("H97").Value = 4159.9
$ws.Range("I97").Value = 616.26086
$ws.Range("J97").Value = 15803.286
$ws.Range("K97").Value = 616.26086
$ws.Range("L97").Value = 15803.286
$ws.Range("M97").Value = -120.26086
$ws.Range("N97").Value = -16795.286
$ws.Range("H107").Value = 946.8276
$ws.Range("I107").Value = 962.2273
$ws.Range("K107").Value = 962.2273
$ws.Range("M107").Value = 957.7727
$ws.Range("H132").Value = 5968027.5
$ws.Range("I132").Value = 2007
$ws.Range("J132").Value = 23866088
$ws.Range("K132").Value = 6021
$ws.Range("L132").Value = 71598264
$ws.Range("M132").Value = -3491
$ws.Range("N132").Value = -71603324

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5879.1904
$ws.Range("I22").Value = 9795.916999999999
$ws.Range("J22").Value = 656.8889
$ws.Range("K22").Value = 9795.916999999999
$ws.Range("L22").Value = 656.8889
$ws.Range("M22").Value = -9500.916999999999
$ws.Range("N22").Value = -1246.8889
$ws.Range("H27").Value = 5879.1904
$ws.Range("I27").Value = 9795.916999999999
$ws.Range("J27").Value = 656.8889
$ws.Range("K27").Value = 9795.916999999999
$ws.Range("L27").Value = 656.8889
$ws.Range("M27").Value = -9688.916999999999
$ws.Range("N27").Value = -870.8889
$ws.Range("H68").Value = 10804957
$ws.Range("I68").Value = 19446082
$ws.Range("J68").Value = 3550.5
$ws.Range("K68").Value = 19446082
$ws.Range("L68").Value = 3550.5
$ws.Range("M68").Value = -19445333
$ws.Range("N68").Value = -5048.5
$ws.Range("H71").Value = 10804957
$ws.Range("I71").Value = 19446082
$ws.Range("J71").Value = 3550.5
$ws.Range("K71").Value = 97230410
$ws.Range("L71").Value = 17752.5
$ws.Range("M71").Value = -97226666
$ws.Range("N71").Value = -25240.5
$ws.Range("H100").Value = 13174823
$ws.Range("I100").Value = 3798.7778
$ws.Range("K100").Value = 3798.7778
$ws.Range("M100").Value = -3257.7778
$ws.Range("H122").Value = 3243.25
$ws.Range("I122").Value = 2815.6086
$ws.Range("K122").Value = 8446.825800000001
$ws.Range("M122").Value = -5996.825800000001
$ws.Range("H132").Value = 4930.8823
$ws.Range("I132").Value = 3675
$ws.Range("J132").Value = 6047.222
$ws.Range("K132").Value = 11025
$ws.Range("L132").Value = 18141.666
$ws.Range("M132").Value = -8495
$ws.Range("N132").Value = -23201.666
$ws.Range("H136").Value = 2047.3966
$ws.Range("I136").Value = 1977.6
$ws.Range("J136").Value = 2202.5
$ws.Range("K136").Value = 5932.799999999999
$ws.Range("L136").Value = 6607.5
$ws.Range("M136").Value = -3382.799999999999
$ws.Range("N136").Value = -11707.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18999.75
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 18999.75
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
$ws.Range("H96").Value = 15282.5
$ws.Range("I96").Value = 8641.429
$ws.Range("K96").Value = 8641.429
$ws.Range("M96").Value = -7268.429
$ws.Range("H100").Value = 990
$ws.Range("I100").Value = 990
$ws.Range("K100").Value = 1980
$ws.Range("M100").Value = -1439
$ws.Range("H107").Value = 3205.2812
$ws.Range("I107").Value = 1208.9
$ws.Range("J107").Value = 6532.5835
$ws.Range("K107").Value = 3626.7
$ws.Range("L107").Value = 19597.7505
$ws.Range("M107").Value = -1706.7
$ws.Range("N107").Value = -23437.7505
$ws.Range("H113").Value = 571.5
$ws.Range("I113").Value = 610.2857
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 1830.8571
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 339.1428999999998
$ws.Range("N113").Value = -5240
$ws.Range("H132").Value = 206049.55
$ws.Range("I132").Value = 1822.0605
$ws.Range("J132").Value = 627268.75
$ws.Range("K132").Value = 5466.181500000001
$ws.Range("L132").Value = 1881806.25
$ws.Range("M132").Value = -2936.181500000001
$ws.Range("N132").Value = -1886866.25
$ws.Range("H136").Value = 130224.98
$ws.Range("I136").Value = 7089.3584
$ws.Range("J136").Value = 719516.9
$ws.Range("K136").Value = 21268.0752
$ws.Range("L136").Value = 2158550.7
$ws.Range("M136").Value = -18718.0752
$ws.Range("N136").Value = -2163650.7
